$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A61").Value = "Andrea Barozzi"
$ws.Range("B61").Value = "Michele Ruele | Avanzi"
$ws.Range("C61").Value = "Luca Frasca | Clitoriders"
$ws.Range("D61").Value = "Antonio Calabrò | Avanzi"
$ws.Range("E61").Value = "Alessandro Fait | RSA United"
$ws.Range("F61").Value = "Thomas Cavagna | MAI UNA GIOIA"
